$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1380.6857
$ws.Range("I15").Value = 1380.6857
$ws.Range("K15").Value = 4142.0571
$ws.Range("M15").Value = -3973.0571
$ws.Range("H33").Value = 3879
$ws.Range("J33").Value = 844.3333
$ws.Range("L33").Value = 844.3333
$ws.Range("N33").Value = -1302.3333
$ws.Range("H103").Value = 1195.25
$ws.Range("I103").Value = 1015.75
$ws.Range("J103").Value = 1374.75
$ws.Range("K103").Value = 3047.25
$ws.Range("L103").Value = 4124.25
$ws.Range("M103").Value = -2461.25
$ws.Range("N103").Value = -5296.25
$ws.Range("H141").Value = 18975.234
$ws.Range("I141").Value = 9402.416999999999
$ws.Range("K141").Value = 28207.251
$ws.Range("M141").Value = -23027.251

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3705282.8
$ws.Range("I2").Value = 5556664.5
$ws.Range("K2").Value = 5556664.5
$ws.Range("M2").Value = -5556551.5
$ws.Range("H32").Value = 3624.027
$ws.Range("I32").Value = 2516.2812
$ws.Range("J32").Value = 10713.6
$ws.Range("K32").Value = 2516.2812
$ws.Range("L32").Value = 10713.6
$ws.Range("M32").Value = -2229.2812
$ws.Range("N32").Value = -11287.6
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("N36").ClearContents()
$ws.Range("H74").Value = 139431.5
$ws.Range("I74").Value = 15732.091
$ws.Range("J74").Value = 592996
$ws.Range("K74").Value = 15732.091
$ws.Range("L74").Value = 592996
$ws.Range("M74").Value = -14858.091
$ws.Range("N74").Value = -594744
$ws.Range("H77").Value = 139431.5
$ws.Range("I77").Value = 15732.091
$ws.Range("J77").Value = 592996
$ws.Range("K77").Value = 78660.455
$ws.Range("L77").Value = 2964980
$ws.Range("M77").Value = -74292.455
$ws.Range("N77").Value = -2973716
$ws.Range("H110").Value = 1463395.1
$ws.Range("I110").Value = 3970268
$ws.Range("J110").Value = 1052.6666
$ws.Range("K110").Value = 3970268
$ws.Range("L110").Value = 1052.6666
$ws.Range("M110").Value = -3968223
$ws.Range("N110").Value = -5142.6666
$ws.Range("H116").Value = 3705282.8
$ws.Range("I116").Value = 5556664.5
$ws.Range("K116").Value = 5556664.5
$ws.Range("M116").Value = -5554370.5
$ws.Range("H122").Value = 633610.5600000001
$ws.Range("I122").Value = 1789.625
$ws.Range("K122").Value = 5368.875
$ws.Range("M122").Value = -2918.875
$ws.Range("H132").Value = 9293.9
$ws.Range("I132").Value = 9956.714
$ws.Range("J132").Value = 7747.3335
$ws.Range("K132").Value = 29870.142
$ws.Range("L132").Value = 23242.0005
$ws.Range("M132").Value = -27340.142
$ws.Range("N132").Value = -28302.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3705282.8
$ws.Range("I3").Value = 5556664.5
$ws.Range("K3").Value = 5556664.5
$ws.Range("M3").Value = -5556550.5
$ws.Range("H22").Value = 1263.9445
$ws.Range("I22").Value = 1183.3334
$ws.Range("K22").Value = 1183.3334
$ws.Range("M22").Value = -1010.3334
$ws.Range("H29").Value = 148583.5
$ws.Range("I29").Value = 221875
$ws.Range("J29").Value = 2000.5
$ws.Range("K29").Value = 221875
$ws.Range("L29").Value = 2000.5
$ws.Range("M29").Value = -221586
$ws.Range("N29").Value = -2578.5
$ws.Range("H94").Value = 2457901.8
$ws.Range("I94").Value = 3367662
$ws.Range("J94").Value = 1548.8
$ws.Range("K94").Value = 3367662
$ws.Range("L94").Value = 1548.8
$ws.Range("M94").Value = -3367211
$ws.Range("N94").Value = -2450.8
$ws.Range("H134").Value = 19112.066
$ws.Range("I134").Value = 21968.7
$ws.Range("K134").Value = 65906.10000000001
$ws.Range("M134").Value = -63371.10000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2655.3044
$ws.Range("I58").Value = 1937.8572
$ws.Range("K58").Value = 1937.8572
$ws.Range("M58").Value = -1734.8572
$ws.Range("H86").Value = 10848.818
$ws.Range("I86").Value = 8970.632
$ws.Range("J86").Value = 13397.786
$ws.Range("K86").Value = 8970.632
$ws.Range("L86").Value = 13397.786
$ws.Range("M86").Value = -7847.632
$ws.Range("N86").Value = -15643.786
$ws.Range("H89").Value = 10848.818
$ws.Range("I89").Value = 8970.632
$ws.Range("J89").Value = 13397.786
$ws.Range("K89").Value = 44853.16
$ws.Range("L89").Value = 66988.92999999999
$ws.Range("M89").Value = -39237.16
$ws.Range("N89").Value = -78220.92999999999
$ws.Range("H132").Value = 52601.75
$ws.Range("I132").Value = 57946.445
$ws.Range("J132").Value = 4499.5
$ws.Range("K132").Value = 173839.335
$ws.Range("L132").Value = 13498.5
$ws.Range("M132").Value = -171309.335
$ws.Range("N132").Value = -18558.5
$ws.Range("H136").Value = 2655.3044
$ws.Range("I136").Value = 1937.8572
$ws.Range("K136").Value = 5813.571599999999
$ws.Range("M136").Value = -3263.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 977.25
$ws.Range("I3").Value = 977.25
$ws.Range("K3").Value = 2931.75
$ws.Range("M3").Value = -2819.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 9459.166999999999
$ws.Range("I2").Value = 1494.7778
$ws.Range("J2").Value = 33352.332
$ws.Range("K2").Value = 1494.7778
$ws.Range("L2").Value = 33352.332
$ws.Range("M2").Value = -1381.7778
$ws.Range("N2").Value = -33578.332
$ws.Range("H97").Value = 993308.5600000001
$ws.Range("I97").Value = 1832744.4
$ws.Range("J97").Value = 1248.091
$ws.Range("K97").Value = 1832744.4
$ws.Range("L97").Value = 1248.091
$ws.Range("M97").Value = -1832248.4
$ws.Range("N97").Value = -2240.091
$ws.Range("H102").Value = 3778344.2
$ws.Range("I102").Value = 5557095.5
$ws.Range("K102").Value = 5557095.5
$ws.Range("M102").Value = -5555473.5
$ws.Range("H132").Value = 12569.044
$ws.Range("I132").Value = 9849.866
$ws.Range("K132").Value = 29549.598
$ws.Range("M132").Value = -27019.598
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1070.8889
$ws.Range("I16").Value = 611
$ws.Range("K16").Value = 611
$ws.Range("M16").Value = -441
$ws.Range("H22").Value = 54147.633
$ws.Range("I22").Value = 111001
$ws.Range("J22").Value = 2979.6
$ws.Range("K22").Value = 111001
$ws.Range("L22").Value = 2979.6
$ws.Range("M22").Value = -110706
$ws.Range("N22").Value = -3569.6
$ws.Range("H27").Value = 54147.633
$ws.Range("I27").Value = 111001
$ws.Range("J27").Value = 2979.6
$ws.Range("K27").Value = 111001
$ws.Range("L27").Value = 2979.6
$ws.Range("M27").Value = -110894
$ws.Range("N27").Value = -3193.6
$ws.Range("H55").Value = 1257.6578
$ws.Range("I55").Value = 1168.5238
$ws.Range("J55").Value = 1367.7646
$ws.Range("K55").Value = 1168.5238
$ws.Range("L55").Value = 1367.7646
$ws.Range("M55").Value = -995.5237999999999
$ws.Range("N55").Value = -1713.7646

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 449.5
$ws.Range("I100").Value = 399
$ws.Range("K100").Value = 798
$ws.Range("M100").Value = -257
$ws.Range("H122").Value = 2727.2
$ws.Range("I122").Value = 2069.5386
$ws.Range("K122").Value = 6208.6158
$ws.Range("M122").Value = -3758.6158
$ws.Range("H132").Value = 22472520
$ws.Range("J132").Value = 1374846.4
$ws.Range("L132").Value = 4124539.2
$ws.Range("N132").Value = -4129599.2
